$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.515.67"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.190.55"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.27%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.189.09"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.60"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.716.97"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.187.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.630.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.91"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.56"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.103"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.86"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.47"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0701"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "419.09"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.998.50"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.11"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.01%  "
$ws.Range("E45").Value = "  -5.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.93"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.61"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("E51").Value = "  -2.67%  "
